$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "HDJ 1"
$ws.Range("A16").Value = "HDJ 2"
$ws.Range("A17").Value = "HDJ 3"
